$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.009416460990906
$ws.Range("B1").Value = 2.797588348388672
$ws.Range("C1").Value = 5.228207111358643
$ws.Range("D1").Value = 2.099475622177124
$ws.Range("E1").Value = 1.17658531665802
